$d = $word.ActiveDocument

# Locate the split point: end of the 2nd occurrence of "Địa chỉ" in the document,
# i.e. right before " theo CCCD: " in the paragraph
#   "Địa chỉ theo CCCD: «Địa_chỉ_theo_CCCD»"
$rng = $d.Content
$rng.Find.MatchWildcards = $false
$n = 0
$splitPos = -1
while ($rng.Find.Execute("Địa chỉ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $n = $n + 1
    if ($n -eq 2) {
        $splitPos = $rng.End
        break
    }
    $rng.Collapse(0)
    $rng.MoveEnd(1, $d.Content.End - $rng.End) | Out-Null
}

# Move the "_GoBack" bookmark from its old position (after "... ngành ") to the
# split point located above.
$old = $d.Bookmarks("_GoBack")
$old.Delete()

$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove " theo CCCD" so "Địa chỉ theo CCCD: " becomes "Địa chỉ: "
$delRange = $d.Range($splitPos, $splitPos + 10)
$delRange.Delete()
